# Generate Report for handback
#
# The 24745f76-... source file has now been handed back (in sync with en-US),
# so its row moves up to be the first data row (row 2) on every sheet, with
# updated "Handed back" status/handback file/handback datetime; the other two
# rows (ffff28ca..., ffffff0b...) shift down by one row each.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "24745f76-4497-4736-9139-e34de63a432a.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"

$ws.Range("A3").Value = "ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"

$ws.Range("A4").Value = "ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md"
$ws.Range("B4").Value = "Handed back: in sync with en-US"
$ws.Range("C4").Value = "Handed back: in sync with en-US"

$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("C5").Value = "Not to be localized"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "24745f76-4497-4736-9139-e34de63a432a.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.zh-cn.xlf"
$ws.Range("D2").Value = "2016-01-28 09:39:41"
$ws.Range("E2").Value = "24745f76-4497-4736-9139-e34de63a432a.md"
$ws.Range("F2").Value = "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.zh-cn.xlf"
$ws.Range("G2").Value = "2016-01-28 09:40:47"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf"
$ws.Range("D3").Value = "2016-01-28 09:37:23"
$ws.Range("E3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md"
$ws.Range("F3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf"
$ws.Range("G3").Value = "2016-01-28 09:38:09"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = "ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md"
$ws.Range("B4").Value = "Handed back: in sync with en-US"
$ws.Range("C4").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf"
$ws.Range("D4").Value = "2016-01-28 09:37:23"
$ws.Range("E4").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md"
$ws.Range("F4").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf"
$ws.Range("G4").Value = "2016-01-28 09:38:09"
$ws.Range("H4").Value = "Include"

$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "24745f76-4497-4736-9139-e34de63a432a.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.de-de.xlf"
$ws.Range("D2").Value = "2016-01-28 09:39:58"
$ws.Range("E2").Value = "24745f76-4497-4736-9139-e34de63a432a.md"
$ws.Range("F2").Value = "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.de-de.xlf"
$ws.Range("G2").Value = "2016-01-28 09:41:10"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf"
$ws.Range("D3").Value = "2016-01-28 09:37:36"
$ws.Range("E3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md"
$ws.Range("F3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf"
$ws.Range("G3").Value = "2016-01-28 09:38:30"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = "ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md"
$ws.Range("B4").Value = "Handed back: in sync with en-US"
$ws.Range("C4").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf"
$ws.Range("D4").Value = "2016-01-28 09:37:36"
$ws.Range("E4").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md"
$ws.Range("F4").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf"
$ws.Range("G4").Value = "2016-01-28 09:38:30"
$ws.Range("H4").Value = "Include"

$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"
